$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold quantities
$ws.Columns("A:A").Insert()

# Header cell for the new Quantity column (match style of other headers, e.g. D1/E1 -> style index 4)
$ws.Range("A1").Value2 = "Quantity"
$ws.Range("D1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Quantity values for each BOM row
$ws.Range("A2").Value2 = 0
$ws.Range("A3").Value2 = 1
$ws.Range("A4").Value2 = 1
$ws.Range("A5").Value2 = 0
$ws.Range("A6").Value2 = 1
$ws.Range("A7").Value2 = 1
$ws.Range("A8").Value2 = 1
$ws.Range("A9").Value2 = 1
$ws.Range("A10").Value2 = 1
$ws.Range("A11").Value2 = 1
$ws.Range("A12").Value2 = 1
$ws.Range("A13").Value2 = 2

# New BOM line: screws, quantity 2, only a note in column E
$ws.Range("E13").Value2 = "#4 3/16"" Panhead Self-Tapping Screws"
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to reflect where editing ended
$ws.Range("A14").Select()
